$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates (header volume number + week-of dates) ---
$ws.Range("A8").Value = "Volume 31   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/26/2024  Through  9/1/2024"

# --- Column width swap (col E <-> col H bestFit widths) ---
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(8).ColumnWidth
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# --- Cells changing from numeric to text ("0" / "***.*"), copy value+style ---
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("N22").Copy($ws.Range("E22"))

# --- Numeric cell updates ---
# Row 14
$ws.Range("M14").Value = -31.578947368421

# Row 15
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("I15").Value = 25
$ws.Range("J15").Value = 29
$ws.Range("K15").Value = -13.793103448275
$ws.Range("L15").Value = 13.636363636363
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = -59.677419354838

# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 166
$ws.Range("J16").Value = 179
$ws.Range("K16").Value = -7.262569832402
$ws.Range("L16").Value = -31.120331950207
$ws.Range("M16").Value = -44.481605351170
$ws.Range("N16").Value = -89.028420356906

# Row 17
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -41.176470588235
$ws.Range("F17").Value = 44
$ws.Range("G17").Value = 48
$ws.Range("H17").Value = -8.333333333333
$ws.Range("I17").Value = 456
$ws.Range("J17").Value = 489
$ws.Range("K17").Value = -6.748466257668
$ws.Range("L17").Value = -6.557377049180
$ws.Range("M17").Value = 6.293706293706
$ws.Range("N17").Value = -42.928660826032

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 28.571428571428
$ws.Range("I18").Value = 94
$ws.Range("J18").Value = 84
$ws.Range("K18").Value = 11.904761904761
$ws.Range("L18").Value = -29.323308270676
$ws.Range("M18").Value = -52.040816326530
$ws.Range("N18").Value = -82.297551789077

# Row 19
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = 27.777777777777
$ws.Range("I19").Value = 243
$ws.Range("J19").Value = 230
$ws.Range("K19").Value = 5.652173913043
$ws.Range("L19").Value = -10
$ws.Range("M19").Value = -16.780821917808
$ws.Range("N19").Value = -60.806451612903

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 105
$ws.Range("J20").Value = 108
$ws.Range("K20").Value = -2.777777777777
$ws.Range("L20").Value = -16.666666666666
$ws.Range("N20").Value = -74.820143884892

# Row 21
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -29.729729729729
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = 7.692307692307
$ws.Range("I21").Value = 1102
$ws.Range("J21").Value = 1125
$ws.Range("K21").Value = -2.044444444444
$ws.Range("L21").Value = -14.771848414539
$ws.Range("M21").Value = -18.970588235294
$ws.Range("N21").Value = -72.374028578591

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 44
$ws.Range("K22").Value = 41.935483870967
$ws.Range("L22").Value = 18.918918918918
$ws.Range("M22").Value = 33.333333333333

# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -62.5
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = -24
$ws.Range("I23").Value = 239
$ws.Range("J23").Value = 245
$ws.Range("K23").Value = -2.448979591836
$ws.Range("L23").Value = -9.469696969696
$ws.Range("M23").Value = 47.530864197530

# Row 24
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 36.363636363636
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 2.702702702702
$ws.Range("I24").Value = 771
$ws.Range("J24").Value = 761
$ws.Range("K24").Value = 1.314060446780
$ws.Range("L24").Value = -0.899742930591
$ws.Range("M24").Value = 16.114457831325

# Row 25
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 36.363636363636
$ws.Range("I25").Value = 197
$ws.Range("J25").Value = 205
$ws.Range("K25").Value = -3.902439024390
$ws.Range("L25").Value = -15.450643776824

# Row 26
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 34
$ws.Range("E26").Value = -38.235294117647
$ws.Range("F26").Value = 84
$ws.Range("G26").Value = 92
$ws.Range("H26").Value = -8.695652173913
$ws.Range("I26").Value = 580
$ws.Range("J26").Value = 664
$ws.Range("K26").Value = -12.650602409638
$ws.Range("L26").Value = -10.631741140215
$ws.Range("M26").Value = -36.193619361936

# Row 27
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("I27").Value = 34
$ws.Range("J27").Value = 34
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -5.555555555555

# Row 28
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 47
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = -21.666666666666
$ws.Range("L28").Value = -18.965517241379

# Row 29
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("M29").Value = -43.75
$ws.Range("N29").Value = -79.820627802690

# Row 30
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("M30").Value = -44.615384615384
$ws.Range("N30").Value = -82.524271844660
